# Class-Updates-Tracker-Java.xlsx update
# piller1-learn about inheritance,methodoverriding,super and final key in java
#
# Summary of change: a new "Day5" row (if-else/switch topic) is inserted
# after the existing "Day4" row, the old combined Day5-Day8 rows shift
# down to Day6-Day9 (renumbered), the old combined OOP row is split into
# a dedicated "Day10" (class/object/overloading/constructor) row and a
# new "Day11" (this keyword) row, and a trailing "Day12" placeholder row
# is appended. E8 gains a real hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 0. Drop existing hyperlinks; we'll re-create all of them in their
#    final positions once the new row has been inserted and the cell
#    values are in place (this engine does not auto-shift hyperlink
#    anchors on row insert).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 1. Insert the new row for "Day5" above the old row 9 ("Day5" in the
#    original numbering, which becomes "Day6"). Everything below
#    (including the old blank-but-styled row 24, and the hyperlinks)
#    shifts down by one row as a side effect.
# ---------------------------------------------------------------------
$ws.Rows.Item(9).Insert()

# ---------------------------------------------------------------------
# 3. Cell values. Columns B..F, rows 3..16 per the target layout.
# ---------------------------------------------------------------------
$ws.Range("E8").Value2 = "https://youtu.be/AmiV9r7ILQI"

$ws.Range("B9").Value2 = "Day5"
$ws.Range("C9").Value2 = "learn about if-else and switch in java"
$ws.Range("D9").Value2 = "https://youtu.be/OTMsKM8OZNQ"

$ws.Range("B10").Value2 = "Day6"
$ws.Range("B11").Value2 = "Day7"
$ws.Range("B12").Value2 = "Day8"

$ws.Range("B13").Value2 = "Day9"
$ws.Range("C13").Value2 = "lean about static keyword in details"
$ws.Range("D13").Value2 = "https://youtu.be/K9Rvor70Aiw"

$ws.Range("B14").Value2 = "Day10"
$ws.Range("C14").Value2 = "oops start and learn about class and objects ,methodoverloading andconsturctor  java"
$ws.Range("D14").Value2 = "https://youtu.be/E98I2pky-hQ"

$ws.Range("B15").Value2 = "Day11"
$ws.Range("C15").Value2 = "this keyword in java"
$ws.Range("D15").Value2 = "https://youtu.be/DzHGyfZH6fA"

$ws.Range("B16").Value2 = "Day12"

# ---------------------------------------------------------------------
# 4. Formatting. Use PasteSpecial(formats) from existing, already
#    correctly-styled cells as templates so the new cells land on the
#    same style families as their neighbours.
# ---------------------------------------------------------------------
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false

# Row 9 (new "Day5" row): B/C like the other label/topic cells, D like a
# bordered hyperlink cell (all four sides), E/F get the "looks like a
# hyperlink but isn't" style used elsewhere in the sheet.
Copy-Format "B10" "B9"
Copy-Format "C10" "C9"
Copy-Format "D11" "D9"

# Rows 10-12 keep the look of the previous Day5/6/7 rows (now shifted).
Copy-Format "D10" "E9"
Copy-Format "D10" "F9"

# Row 16 (new trailing "Day12" placeholder row): match the blank rows.
Copy-Format "B17" "B16"
Copy-Format "C17" "C16"
Copy-Format "D17" "D16"
Copy-Format "E17" "E16"
Copy-Format "F17" "F16"

$excel.CutCopyMode = $false

# E9/F9 direct styling: underlined black Yu Gothic font, boxed border
# (no bottom rule), centered -- matches the "looks like a link" cells
# used elsewhere (e.g. old E9/F9, E10/F10) but without an actual
# hyperlink relationship.
foreach ($addr in @("E9", "F9")) {
    $r = $ws.Range($addr)
    $r.Font.Name = "游ゴシック"
    $r.Font.Underline = 2
    $r.Font.Color = 0
    $r.HorizontalAlignment = -4108
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(7).Color = 0
    $r.Borders.Item(8).LineStyle = 1
    $r.Borders.Item(8).Color = 0
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(10).Color = 0
    $r.Borders.Item(9).LineStyle = 0
}

# ---------------------------------------------------------------------
# 5. Hyperlinks. Re-create every hyperlink at its final location. Excel
#    auto-applies the "looks like a hyperlink" styling to these cells.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D10"), "https://youtu.be/30xvczNMAUs") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "https://youtu.be/ITr0uX2Ez1o") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E10"), "https://youtu.be/0r1SfRoLuzU") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E11"), "https://youtu.be/qqRDHzPli3o") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://youtu.be/aSGa3S2-sQ0") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D12"), "https://youtu.be/ZnfOBRS7KCg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), "https://youtu.be/E98I2pky-hQ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), "https://youtu.be/DzHGyfZH6fA") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E8"), "https://youtu.be/AmiV9r7ILQI") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "https://youtu.be/OTMsKM8OZNQ") | Out-Null

# ---------------------------------------------------------------------
# 6. Row heights -- new rows should look like their neighbours.
# ---------------------------------------------------------------------
$ws.Rows.Item(9).RowHeight = $ws.Rows.Item(10).RowHeight
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(17).RowHeight

# ---------------------------------------------------------------------
# 7. Sheet view tweaks observed in the diff: selection moved to C18 and
#    the frozen/anchor top-left cell reset to default.
# ---------------------------------------------------------------------
$ws.Range("C18").Select() | Out-Null
